# Update the player-name spellings in the groups list (TournamentPlayer
# model rework) and move the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C corrections
$ws.Range("C5").Value = "Sung-jae Im"
$ws.Range("C7").Value = "Byeong-Hun An"
$ws.Range("C8").Value = "Rafael Cabrera Bello"

# Column D corrections
$ws.Range("D14").Value = "Austin Connelly"
$ws.Range("D25").Value = "Dong-Kyu Jang"
$ws.Range("D33").Value = "Brandon Wu (a)"
$ws.Range("D35").Value = "Inn-choon Hwang"
$ws.Range("D37").Value = "Takumi Kanaya (a)"
$ws.Range("D41").Value = "Curtis Knipes (a)"
$ws.Range("D50").Value = "Li Haotong"
$ws.Range("D55").Value = "Mike Lorenzo-Vera"
$ws.Range("D62").Value = "Alexander Noren"
$ws.Range("D67").Value = "Dimitrios Papadatos"
$ws.Range("D68").Value = "Sang-hyun Park"
$ws.Range("D76").Value = "Matthias Schmid (a)"
$ws.Range("D85").Value = "James Sugrue (a)"
$ws.Range("D88").Value = "Thomas Thurloway (a)"

# Move the active selection to match the author's last cursor position
$ws.Range("G18").Select()
